# [2023. 01. 09] 1차 test 버전 - 보완(10:25)
#
# Updates a handful of "Text" / "Rectangle" shapes on slide 1 of the
# weekly-summary deck: status dates turn into "[보류]" placeholders,
# a couple of task descriptions are rewritten, two status rectangles
# go from green back to white, and several trailing detail fields
# (ticket id / SQL note / date / user id) get refreshed values.

function Set-ShapeLineText {
    # $PartIndex is the 0-based index into Text.Split([char]11) — i.e. the
    # same split PowerPoint uses internally where a soft line break
    # (<a:br/>) shows up as a vertical-tab (chr 11) character. Shapes whose
    # body starts with a leading <a:br/> have an empty part[0], so their
    # first visible line is part[1], and so on.
    param(
        $Shape,
        [int]$PartIndex,
        [string]$NewText
    )

    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $vt = [char]11
    $parts = $full.Split($vt)

    $offset = 1
    for ($i = 0; $i -lt $PartIndex; $i++) {
        $offset += $parts[$i].Length + 1
    }
    $len = $parts[$PartIndex].Length

    $sub = $tr.Characters($offset, $len)
    $sub.Text = $NewText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 8: "12/22" -> "[보류]"   (body is <br/><r>12/22</r><br/>)
Set-ShapeLineText $s.Shapes.Item(8) 1 "[보류]"

# Shape 9: " [E-BIZ] 거래처 정보화면 임원 별도처리" -> " [e-Biz] 매출원장(ZSDR5370)보완"
Set-ShapeLineText $s.Shapes.Item(9) 1 " [e-Biz] 매출원장(ZSDR5370)보완"

# Shape 12: two lines rewritten
Set-ShapeLineText $s.Shapes.Item(12) 1 " [e-Biz] 매출원장(ZSDR5370)보완"
Set-ShapeLineText $s.Shapes.Item(12) 2 " [e-Biz/FLBIZ] 감사 대비 어플리케이션 로그 기능"

# Shape 27: "12/02" -> "[보류]", "12/13" -> "01/02"
Set-ShapeLineText $s.Shapes.Item(27) 1 "[보류]"
Set-ShapeLineText $s.Shapes.Item(27) 2 "01/02"

# Shape 29: " 1" -> " [RMS] 시스템 개발 및 보완", " 2" -> " [RMS] 1차 테스트"
Set-ShapeLineText $s.Shapes.Item(29) 1 " [RMS] 시스템 개발 및 보완"
Set-ShapeLineText $s.Shapes.Item(29) 2 " [RMS] 1차 테스트"

# Shape 31: "미진행" -> "진행중"   (single run, no surrounding <br/>)
Set-ShapeLineText $s.Shapes.Item(31) 0 "진행중"

# Shape 32: "[보류]" -> "90%" (first line), second "[보류]" line unchanged
Set-ShapeLineText $s.Shapes.Item(32) 1 "90%"

# Shape 35: " t" -> " [RMS] 1차 테스트", " s" -> " [RMS] 2차 테스트"
Set-ShapeLineText $s.Shapes.Item(35) 1 " [RMS] 1차 테스트"
Set-ShapeLineText $s.Shapes.Item(35) 2 " [RMS] 2차 테스트"

# Shapes 37 & 38: status rectangles go from green back to white
$s.Shapes.Item(37).Fill.ForeColor.RGB = 0xFFFFFF
$s.Shapes.Item(38).Fill.ForeColor.RGB = 0xFFFFFF

# Shape 40: "Z11-22-0019" -> "4718"   (single run, no surrounding <br/>)
Set-ShapeLineText $s.Shapes.Item(40) 0 "4718"

# Shape 41: "Field contents changed: I_PSTYP -> 0" -> SQL note
Set-ShapeLineText $s.Shapes.Item(41) 0 'update summary set sign="승인" where sum_id=1'

# Shape 42: "2022.12.08" -> "2023-01-05"
Set-ShapeLineText $s.Shapes.Item(42) 0 "2023-01-05"

# Shape 50: "SROH" -> "jelee01"
Set-ShapeLineText $s.Shapes.Item(50) 0 "jelee01"
